$wb = $excel.ActiveWorkbook

# Update the status text "Ready for handoff" -> "In Translation"
# on every sheet that contains it (Overview, zh-cn, de-de).
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # Cast to [string] explicitly: some cell values (e.g. the literal
        # text "True") come back as native booleans, and PowerShell's -eq
        # would otherwise coerce the right-hand string into a boolean too
        # (any non-empty string -> $true), causing false-positive matches.
        $text = [string]$cell.Value2
        if ($text -ceq "Ready for handoff") {
            $cell.Value = "In Translation"
        }
    }
}

# Narrow the "Status" column(s) on all three sheets.
# (stored width 17.2159881591797 -> 13.4101845877511; ColumnWidth snaps to the
# nearest 1/6-character grid, so 12.5 lands on the closest achievable value.)
$newWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E1:F1").ColumnWidth = $newWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C1").ColumnWidth = $newWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C1").ColumnWidth = $newWidth
